$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Recorded By" (column G) values to reflect reordering of names/emails
$updates = @{
    2 = "backup@backdoor.com, system, System"
    4 = "backup@backdoor.com, System"
    5 = "backup@backdoor.com, System"
    7 = "admin@admin.com, System"
    10 = "dnasr281@gmail.com, System"
    11 = "System, dnasr281@gmail.com"
    17 = "System, dnasr281@gmail.com"
    19 = "dnasr281@gmail.com, System"
    29 = "backup@backdoor.com, system, System"
    31 = "backup@backdoor.com, System"
    32 = "backup@backdoor.com, System"
    34 = "admin@admin.com, System"
    37 = "dnasr281@gmail.com, System"
    38 = "System, dnasr281@gmail.com"
    44 = "System, dnasr281@gmail.com"
    46 = "dnasr281@gmail.com, System"
    56 = "backup@backdoor.com, system, System"
    58 = "backup@backdoor.com, System"
    59 = "backup@backdoor.com, System"
    61 = "admin@admin.com, System"
    64 = "dnasr281@gmail.com, System"
    65 = "System, dnasr281@gmail.com"
    71 = "System, dnasr281@gmail.com"
    73 = "dnasr281@gmail.com, System"
    83 = "backup@backdoor.com, System"
    84 = "backup@backdoor.com, System"
    85 = "backup@backdoor.com, System"
    86 = "dnasr281@gmail.com, System"
    90 = "admin@admin.com, dnasr281@gmail.com"
    96 = "System, dnasr281@gmail.com"
    109 = "backup@backdoor.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    112 = "dnasr281@gmail.com, System"
    116 = "admin@admin.com, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    135 = "backup@backdoor.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    138 = "dnasr281@gmail.com, System"
    142 = "admin@admin.com, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $updates[$row]
}

